# "Add files via upload" -- adds two labeled cells (F14/G14) to the
# sheet and repositions the two existing charts further down/right on
# the sheet (as if the author dragged them to make room for the new
# row-14 content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New cell content (row 14 gains two more labels) -------------
$ws.Range("F14").Value = "On diff"
$ws.Range("G14").Value = "comps"

# --- 2. Reposition the two chart objects -----------------------------
# ChartObjects(1) = "Діаграма 3"  (the "1..6 threads" line chart)
# ChartObjects(2) = "Best tests"  (the "second"/"main" line chart)
$cos = $ws.ChartObjects()

$co1 = $cos.Item(1)
$co1.Left = 469.71102362204726
$co1.Top = 202.41889763779528
$co1.Width = 433.04360236220475
$co1.Height = 160.748031496063

$co2 = $cos.Item(2)
$co2.Left = 84.37450787401575
$co2.Top = 217.84409448818897
$co2.Width = 373.89744094488185
$co2.Height = 128.90078740157483

# --- 3. Update the active selection ----------------------------------
$ws.Range("G19").Select() | Out-Null
